$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.461.27'
$ws.Cells.Item(2, 5).Value = '  +0.49%  '
$ws.Cells.Item(3, 4).Value = '3.089.62'
$ws.Cells.Item(3, 5).Value = '  +1.20%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '516.83'
$ws.Cells.Item(5, 5).Value = '  +0.33%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '141.58'
$ws.Cells.Item(6, 5).Value = '  +0.07%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 5).Value = '  -1.10%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '7.28'
$ws.Cells.Item(9, 5).Value = '  +0.75%  '
$ws.Cells.Item(10, 5).Value = '  -0.91%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.374'
$ws.Cells.Item(11, 5).Value = '  -1.10%  '
$ws.Cells.Item(12, 4).Value = '3.615.98'
$ws.Cells.Item(12, 5).Value = '  +1.26%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '25.65'
$ws.Cells.Item(14, 5).Value = '  -5.13%  '
$ws.Cells.Item(15, 5).Value = '  -1.73%  '
$ws.Cells.Item(16, 4).Value = '57.574.57'
$ws.Cells.Item(16, 5).Value = '  +0.78%  '
$ws.Cells.Item(17, 4).Value = '3.088.18'
$ws.Cells.Item(17, 5).Value = '  +1.30%  '
$ws.Cells.Item(18, 5).Value = '  -0.69%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.24'
$ws.Cells.Item(19, 5).Value = '  -1.40%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.15'
$ws.Cells.Item(20, 5).Value = '  +0.00%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '334.53'
$ws.Cells.Item(21, 5).Value = '  +0.93%  '
$ws.Cells.Item(22, 5).Value = '  +0.15%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.501'
$ws.Cells.Item(23, 5).Value = '  -1.30%  '
$ws.Cells.Item(24, 5).Value = '  +0.55%  '
$ws.Cells.Item(25, 5).Value = '  +3.15%  '
$ws.Cells.Item(26, 5).Value = '  +0.15%  '
$ws.Cells.Item(27, 4).Value = '0.0₃0910'
$ws.Cells.Item(27, 5).Value = '  +1.58%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.39'
$ws.Cells.Item(28, 5).Value = '  -5.61%  '
$ws.Cells.Item(29, 5).Value = '  -0.83%  '
$ws.Cells.Item(30, 5).Value = '  +0.33%  '
$ws.Cells.Item(31, 5).Value = '  +0.17%  '
$ws.Cells.Item(32, 5).Value = '  -3.64%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '154.79'
$ws.Cells.Item(33, 5).Value = '  +2.76%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '28.00'
$ws.Cells.Item(34, 5).Value = '  +10.34%  '
$ws.Cells.Item(35, 5).Value = '  -3.72%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.88'
$ws.Cells.Item(36, 5).Value = '  -1.30%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.26'
$ws.Cells.Item(37, 5).Value = '  -1.41%  '
$ws.Cells.Item(38, 5).Value = '  -0.52%  '
$ws.Cells.Item(39, 4).Value = '3.130.12'
$ws.Cells.Item(39, 5).Value = '  +1.48%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '36.76'
$ws.Cells.Item(40, 5).Value = '  +0.05%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.87'
$ws.Cells.Item(41, 5).Value = '  -0.85%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.671'
$ws.Cells.Item(42, 5).Value = '  +1.19%  '
$ws.Cells.Item(43, 5).Value = '  +0.04%  '
$ws.Cells.Item(44, 4).Value = '2.291.69'
$ws.Cells.Item(44, 5).Value = '  +3.88%  '
$ws.Cells.Item(45, 5).Value = '  +5.85%  '
$ws.Cells.Item(46, 5).Value = '  -1.43%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.939'
$ws.Cells.Item(47, 5).Value = '  -1.30%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '20.02'
$ws.Cells.Item(48, 5).Value = '  -0.88%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '5.89'
$ws.Cells.Item(49, 5).Value = '  -3.41%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '253.41'
$ws.Cells.Item(50, 5).Value = '  +5.91%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0877'
$ws.Cells.Item(51, 5).Value = '  +0.90%  '
